# Corrected Technologies and removed values too close to an insert table in the stoch file
# Insert a new row after row 7 (pushing rows 8-15 down to 9-16) and fill the
# new row 8 with a copy of row 7's content.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows.Item(8).EntireRow.Insert()
$ws.Range("A7:K7").Copy($ws.Range("A8:K8"))
